# Apply the updates described by the commit "updated task used in testing":
#  - D2: 4 -> 5
#  - F2: 2 -> 3
#  - H2: 36 -> 46
#  - Active selection moves from D4 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

$ws.Range("D2").Select()
